$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dces")

# Insert two blank rows above row 2, pushing the existing data (rows 2-10) down to rows 4-12
$ws.Rows("2:3").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Add the new row of data at row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 4).Value = 56
$ws.Cells.Item(13, 10).Value = "abc"

# Update the selection to match the target state
$ws.Range("A4:J13").Select()
